$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: "Marking" row
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12: "Total" row
$ws.Range("B12").Value = 96
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "88 / 112"
